# Actualización automática 2025-09-18 09:12:30
# Updates sales figures across the three worksheets to reflect newly
# reported data for "GUERRERO FAREZ FABIAN MAURICIO".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("E4").Value = 138.26
$wsGrupo.Range("M4").Value = 1654.13
$wsGrupo.Range("N4").Value = 234.38

$wsGrupo.Range("L5").Value = 2792.5

$wsGrupo.Range("M7").Value = 1240.62

$wsGrupo.Range("O8").Value = 529.1799999999999

$wsGrupo.Range("L12").Value = 582.53

$wsGrupo.Range("M13").Value = -1790.7

$wsGrupo.Range("D21").Value = 915.84
$wsGrupo.Range("I21").Value = 389.7
$wsGrupo.Range("L21").Value = 855.36
$wsGrupo.Range("M21").Value = 3254.7

$wsGrupo.Range("M25").Value = -1054.31

$wsGrupo.Range("M26").Value = 225.89

$wsGrupo.Range("M38").Value = 297.16

$wsGrupo.Range("E51").Value = 111.3
$wsGrupo.Range("M51").Value = 221.62

$wsGrupo.Range("L54").Value = "7 de 52"
$wsGrupo.Range("N54").Value = "1 de 52"
$wsGrupo.Range("O54").Value = "1 de 52"

# ---------------------------------------------------------------------
# Sheet 2: VENTA MENSUAL
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F4").Value = 3956.81
$wsMensual.Range("F5").Value = 5483.48
$wsMensual.Range("F7").Value = 1240.62
$wsMensual.Range("F8").Value = 529.1799999999999
$wsMensual.Range("F12").Value = 582.53
$wsMensual.Range("F13").Value = -1790.7
$wsMensual.Range("F21").Value = 5480.42
$wsMensual.Range("F25").Value = -1054.31
$wsMensual.Range("F26").Value = 683.8099999999999
$wsMensual.Range("F38").Value = 297.16
$wsMensual.Range("F53").Value = 575.2
$wsMensual.Range("F54").Value = 575.2
$wsMensual.Range("F58").Value = 34176.34

# ---------------------------------------------------------------------
# Sheet 3: CUMPLIMIENTO MENSUAL
# ---------------------------------------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumplimiento.Range("D3").Value = 2222.2
$wsCumplimiento.Range("E3").Value = 15446.9470988183
$wsCumplimiento.Range("F3").Value = 0.1257672477099146

$wsCumplimiento.Range("D4").Value = 314.38
$wsCumplimiento.Range("E4").Value = 728.84288526528
$wsCumplimiento.Range("F4").Value = 0.3013545853339449

$wsCumplimiento.Range("D7").Value = 2428.2
$wsCumplimiento.Range("E7").Value = -1541.488983712426
$wsCumplimiento.Range("F7").Value = 2.738434456545082

$wsCumplimiento.Range("D11").Value = 8329.26
$wsCumplimiento.Range("E11").Value = 9502.1543984654
$wsCumplimiento.Range("F11").Value = 0.4671115714027054

$wsCumplimiento.Range("D12").Value = 15538.08
$wsCumplimiento.Range("E12").Value = 46325.6403947566
$wsCumplimiento.Range("F12").Value = 0.2511662716184939

$wsCumplimiento.Range("D13").Value = 234.38
$wsCumplimiento.Range("E13").Value = 206.273177778119
$wsCumplimiento.Range("F13").Value = 0.5318922268569609

$wsCumplimiento.Range("D14").Value = 529.1799999999999
$wsCumplimiento.Range("E14").Value = 7308.13410570622
$wsCumplimiento.Range("F14").Value = 0.06752058075798092

$wsCumplimiento.Range("D15").Value = 33025.58
$wsCumplimiento.Range("E15").Value = 89029.25551083435
$wsCumplimiento.Range("F15").Value = 0.2705798575023965
